$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and its "through" date label
$ws.Name = "Through 2022-03-24"
$ws.Range("A4").Value = "March (through 03-24)"

# Update March row (row 4) values
$ws.Range("C4").Value = 35
$ws.Range("D4").Value = 45
$ws.Range("F4").Value = 25
$ws.Range("G4").Value = 46
$ws.Range("H4").Value = 64
$ws.Range("I4").Value = 99

# Update Total row (row 5) values
$ws.Range("C5").Value = 122
$ws.Range("D5").Value = 176
$ws.Range("F5").Value = 104
$ws.Range("G5").Value = 187
$ws.Range("H5").Value = 406
$ws.Range("I5").Value = 399
